$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H88").Value = 6180
$ws.Range("I88").Value = 732
$ws.Range("J88").Value = 9293.143
$ws.Range("K88").Value = 732
$ws.Range("L88").Value = 9293.143
$ws.Range("M88").Value = -326
$ws.Range("N88").Value = -10105.143
$ws.Range("H91").Value = 6180
$ws.Range("I91").Value = 732
$ws.Range("J91").Value = 9293.143
$ws.Range("K91").Value = 732
$ws.Range("L91").Value = 9293.143
$ws.Range("M91").Value = 672
$ws.Range("N91").Value = -12101.143
$ws.Range("H100").Value = 1414.8572
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H106").Value = 77500
$ws.Range("I106").Value = 77500
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 77500
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -76869
$ws.Range("N106").ClearContents()
$ws.Range("H111").Value = 1999
$ws.Range("I111").Value = 1999
$ws.Range("K111").Value = 5997
$ws.Range("M111").Value = -2930
$ws.Range("H112").Value = 2516.2727
$ws.Range("J112").Value = 2426.4285
$ws.Range("L112").Value = 7279.2855
$ws.Range("N112").Value = -9495.2855
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H132").Value = 823.6667
$ws.Range("I132").Value = 908.4
$ws.Range("K132").Value = 2725.2
$ws.Range("M132").Value = -195.1999999999998
$ws.Range("H137").Value = 2142.4285
$ws.Range("I137").Value = 1999.6666
$ws.Range("J137").Value = 2249.5
$ws.Range("K137").Value = 5998.9998
$ws.Range("L137").Value = 6748.5
$ws.Range("M137").Value = -3448.9998
$ws.Range("N137").Value = -11848.5
$ws.Range("H138").Value = 2144.525
$ws.Range("I138").Value = 1917.2333
$ws.Range("K138").Value = 5751.699900000001
$ws.Range("M138").Value = -611.6999000000005
$ws.Range("H141").Value = 1484
$ws.Range("I141").Value = 1500.6522
$ws.Range("J141").Value = 1292.5
$ws.Range("K141").Value = 4501.9566
$ws.Range("L141").Value = 3877.5
$ws.Range("M141").Value = 678.0434000000005
$ws.Range("N141").Value = -14237.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 496.66666
$ws.Range("I88").Value = 420
$ws.Range("K88").Value = 420
$ws.Range("M88").Value = -14
$ws.Range("H91").Value = 496.66666
$ws.Range("I91").Value = 420
$ws.Range("K91").Value = 420
$ws.Range("M91").Value = 984
$ws.Range("H132").Value = 2115
$ws.Range("I132").Value = 2115
$ws.Range("K132").Value = 6345
$ws.Range("M132").Value = -3815

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 542.2308
$ws.Range("I22").Value = 559
$ws.Range("K22").Value = 559
$ws.Range("M22").Value = -386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5639.7
$ws.Range("I16").Value = 5730.6
$ws.Range("J16").Value = 5548.8
$ws.Range("K16").Value = 5730.6
$ws.Range("L16").Value = 5548.8
$ws.Range("M16").Value = -5443.6
$ws.Range("N16").Value = -6122.8
$ws.Range("H26").Value = 4419
$ws.Range("J26").Value = 4419
$ws.Range("L26").Value = 4419
$ws.Range("N26").Value = -4993
$ws.Range("H29").Value = 7291
$ws.Range("J29").Value = 7291
$ws.Range("L29").Value = 7291
$ws.Range("N29").Value = -7877
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 5639.7
$ws.Range("I113").Value = 5730.6
$ws.Range("J113").Value = 5548.8
$ws.Range("K113").Value = 5730.6
$ws.Range("L113").Value = 5548.8
$ws.Range("M113").Value = -3560.6
$ws.Range("N113").Value = -9888.799999999999
$ws.Range("H132").Value = 1935.8125
$ws.Range("I132").Value = 1988.8572
$ws.Range("J132").Value = 1564.5
$ws.Range("K132").Value = 5966.571599999999
$ws.Range("L132").Value = 4693.5
$ws.Range("M132").Value = -3436.571599999999
$ws.Range("N132").Value = -9753.5
$ws.Range("H134").Value = 1175.4
$ws.Range("I134").Value = 1175.4
$ws.Range("K134").Value = 3526.2
$ws.Range("M134").Value = -991.2000000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 9000
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 11000
$ws.Range("K64").Value = 9000
$ws.Range("L64").Value = 33000
$ws.Range("M64").Value = -8730
$ws.Range("N64").Value = -33540
$ws.Range("H67").Value = 9000
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 11000
$ws.Range("K67").Value = 9000
$ws.Range("L67").Value = 33000
$ws.Range("M67").Value = -8064
$ws.Range("N67").Value = -34872
$ws.Range("H98").Value = 2101.3333
$ws.Range("J98").Value = 4003.6667
$ws.Range("L98").Value = 12011.0001
$ws.Range("N98").Value = -15007.0001
$ws.Range("H104").Value = 2683
$ws.Range("I104").Value = 2683
$ws.Range("K104").Value = 8049
$ws.Range("M104").Value = -5428
$ws.Range("H121").Value = 15744.9
$ws.Range("I121").Value = 25624
$ws.Range("J121").Value = 5865.8
$ws.Range("K121").Value = 76872
$ws.Range("L121").Value = 17597.4
$ws.Range("M121").Value = -75562
$ws.Range("N121").Value = -20217.4
$ws.Range("H137").Value = 4623.5
$ws.Range("I137").Value = 2300
$ws.Range("J137").Value = 5398
$ws.Range("K137").Value = 6900
$ws.Range("L137").Value = 16194
$ws.Range("M137").Value = -1800
$ws.Range("N137").Value = -26394

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6690.6
$ws.Range("I70").Value = 6358.9
$ws.Range("J70").Value = 7354
$ws.Range("K70").Value = 6358.9
$ws.Range("L70").Value = 7354
$ws.Range("M70").Value = -6088.9
$ws.Range("N70").Value = -7894
$ws.Range("H73").Value = 6690.6
$ws.Range("I73").Value = 6358.9
$ws.Range("J73").Value = 7354
$ws.Range("K73").Value = 6358.9
$ws.Range("L73").Value = 7354
$ws.Range("M73").Value = -5422.9
$ws.Range("N73").Value = -9226
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -2996
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 5000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -14984
$ws.Range("H107").Value = 3593.5557
$ws.Range("I107").Value = 2499.6667
$ws.Range("J107").Value = 4140.5
$ws.Range("K107").Value = 2499.6667
$ws.Range("L107").Value = 4140.5
$ws.Range("M107").Value = -579.6667000000002
$ws.Range("N107").Value = -7980.5
$ws.Range("H126").Value = 2717.7144
$ws.Range("I126").Value = 2552.75
$ws.Range("J126").Value = 2937.6667
$ws.Range("K126").Value = 7658.25
$ws.Range("L126").Value = 8813.000100000001
$ws.Range("M126").Value = -5188.25
$ws.Range("N126").Value = -13753.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 76999.7
$ws.Range("I2").Value = 97142.86
$ws.Range("K2").Value = 97142.86
$ws.Range("M2").Value = -97030.86
$ws.Range("H7").Value = 5914.6924
$ws.Range("I7").Value = 3026.6365
$ws.Range("K7").Value = 3026.6365
$ws.Range("M7").Value = -2914.6365
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H101").Value = 19180
$ws.Range("J101").Value = 19180
$ws.Range("L101").Value = 19180
$ws.Range("N101").Value = -25670
$ws.Range("H126").Value = 5914.6924
$ws.Range("I126").Value = 3026.6365
$ws.Range("K126").Value = 9079.9095
$ws.Range("M126").Value = -6609.9095
$ws.Range("H136").Value = 2864.5
$ws.Range("I136").Value = 2416.4285
$ws.Range("K136").Value = 7249.2855
$ws.Range("M136").Value = -4699.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 541249.75
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 54999.668
$ws.Range("K2").Value = 2000000
$ws.Range("L2").Value = 54999.668
$ws.Range("M2").Value = -1999888
$ws.Range("N2").Value = -55223.668
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H132").Value = 2596.5454
$ws.Range("I132").Value = 2848.4119
$ws.Range("J132").Value = 1740.2
$ws.Range("K132").Value = 8545.235700000001
$ws.Range("L132").Value = 5220.6
$ws.Range("M132").Value = -6015.235700000001
$ws.Range("N132").Value = -10280.6
$ws.Range("H136").Value = 2585.6155
$ws.Range("I136").Value = 2772.7917
$ws.Range("J136").Value = 339.5
$ws.Range("K136").Value = 8318.375100000001
$ws.Range("L136").Value = 1018.5
$ws.Range("M136").Value = -5768.375100000001
$ws.Range("N136").Value = -6118.5
